# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update "Datos actualizados" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Julio de 2020 a las 09:46"

# 2. Paraguay/Hungria: rows 102-103 swap order (Hungria now listed before Paraguay)
#    and Hungria's stats get updated; Paraguay keeps its previous stats.
$ws.Range("A102").Value = "Hungria"
$ws.Range("B102").Value = 4448
$ws.Range("C102").Value = 13
$ws.Range("D102").Value = 3329
$ws.Range("E102").Value = 523
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 596

$ws.Range("A103").Value = "Paraguay"
$ws.Range("B103").Value = 4444
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 2794
$ws.Range("E103").Value = 1609
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 41

# 3. Islas Malvinas/Groenlandia: rows 210-211 swap order (Islas Malvinas now before Groenlandia)
#    Their stats are identical, so only the country names are swapped.
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"

# 4. Update statistics for several other countries (rows identified by row number)

# Row 7 - Rusia
$ws.Range("B7").Value = 818120
$ws.Range("C7").Value = 5635
$ws.Range("D7").Value = 603329
$ws.Range("E7").Value = 201437
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 85
$ws.Range("H7").Value = 13354

# Row 38 - Ucrania
$ws.Range("B38").Value = 65656
$ws.Range("C38").Value = 807
$ws.Range("D38").Value = 36112
$ws.Range("E38").Value = 27928
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 11
$ws.Range("H38").Value = 1616

# Row 53 - Armenia
$ws.Range("B53").Value = 37390
$ws.Range("C53").Value = 73
$ws.Range("D53").Value = 26665
$ws.Range("E53").Value = 10014
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 711

# Row 124 - Eslovaquia
$ws.Range("B124").Value = 2181
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 1616
$ws.Range("E124").Value = 537
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 28

# Row 126 - Estonia
$ws.Range("D126").Value = 1923
$ws.Range("E126").Value = 42

# Row 127 - Lituania
$ws.Range("B127").Value = 2019
$ws.Range("C127").Value = 11
$ws.Range("D127").Value = 1620
$ws.Range("E127").Value = 319
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 80

# Row 162 - Taiwan
$ws.Range("B162").Value = 462
$ws.Range("C162").Value = 4
$ws.Range("D162").Value = 440
$ws.Range("E162").Value = 15
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 7
